$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 5270, 45747),
    @(3, 5220, 45747.01041666666),
    @(4, 5180, 45747.02083333334),
    @(5, 5140, 45747.03125),
    @(6, 5100, 45747.04166666666),
    @(7, 5080, 45747.05208333334),
    @(8, 5070, 45747.0625),
    @(9, 5060, 45747.07291666666),
    @(10, 5050, 45747.08333333334),
    @(11, 5050, 45747.09375),
    @(12, 5050, 45747.10416666666),
    @(13, 5060, 45747.11458333334),
    @(14, 5080, 45747.125),
    @(15, 5110, 45747.13541666666),
    @(16, 5130, 45747.14583333334),
    @(17, 5160, 45747.15625),
    @(18, 5220, 45747.16666666666),
    @(19, 5280, 45747.17708333334),
    @(20, 5360, 45747.1875),
    @(21, 5470, 45747.19791666666),
    @(22, 5620, 45747.20833333334),
    @(23, 5760, 45747.21875),
    @(24, 5920, 45747.22916666666),
    @(25, 6100, 45747.23958333334),
    @(26, 6340, 45747.25),
    @(27, 6540, 45747.26041666666),
    @(28, 6720, 45747.27083333334),
    @(29, 6880, 45747.28125),
    @(30, 7040, 45747.29166666666),
    @(31, 7180, 45747.30208333334),
    @(32, 7260, 45747.3125),
    @(33, 7320, 45747.32291666666),
    @(34, 7350, 45747.33333333334),
    @(35, 7360, 45747.34375),
    @(36, 7350, 45747.35416666666),
    @(37, 7340, 45747.36458333334),
    @(38, 7300, 45747.375),
    @(39, 7260, 45747.38541666666),
    @(40, 7230, 45747.39583333334),
    @(41, 7170, 45747.40625),
    @(42, 7100, 45747.41666666666),
    @(43, 7050, 45747.42708333334),
    @(44, 6990, 45747.4375),
    @(45, 6930, 45747.44791666666),
    @(46, 6860, 45747.45833333334),
    @(47, 6790, 45747.46875),
    @(48, 6710, 45747.47916666666),
    @(49, 6640, 45747.48958333334),
    @(50, 6580, 45747.5),
    @(51, 6550, 45747.51041666666),
    @(52, 6510, 45747.52083333334),
    @(53, 6500, 45747.53125),
    @(54, 6490, 45747.54166666666),
    @(55, 6480, 45747.55208333334),
    @(56, 6480, 45747.5625),
    @(57, 6490, 45747.57291666666),
    @(58, 6510, 45747.58333333334),
    @(59, 6520, 45747.59375),
    @(60, 6540, 45747.60416666666),
    @(61, 6560, 45747.61458333334),
    @(62, 6580, 45747.625),
    @(63, 6610, 45747.63541666666),
    @(64, 6660, 45747.64583333334),
    @(65, 6710, 45747.65625),
    @(66, 6770, 45747.66666666666),
    @(67, 6830, 45747.67708333334),
    @(68, 6900, 45747.6875),
    @(69, 6960, 45747.69791666666),
    @(70, 7030, 45747.70833333334),
    @(71, 7090, 45747.71875),
    @(72, 7160, 45747.72916666666),
    @(73, 7240, 45747.73958333334),
    @(74, 7310, 45747.75),
    @(75, 7380, 45747.76041666666),
    @(76, 7460, 45747.77083333334),
    @(77, 7540, 45747.78125),
    @(78, 7640, 45747.79166666666),
    @(79, 7700, 45747.80208333334),
    @(80, 7710, 45747.8125),
    @(81, 7710, 45747.82291666666),
    @(82, 7660, 45747.83333333334),
    @(83, 7570, 45747.84375),
    @(84, 7450, 45747.85416666666),
    @(85, 7300, 45747.86458333334),
    @(86, 7140, 45747.875),
    @(87, 7000, 45747.88541666666),
    @(88, 6830, 45747.89583333334),
    @(89, 6660, 45747.90625),
    @(90, 6500, 45747.91666666666),
    @(91, 6360, 45747.92708333334),
    @(92, 6240, 45747.9375),
    @(93, 6120, 45747.94791666666),
    @(94, 5900, 45747.95833333334),
    @(95, 5850, 45747.96875),
    @(96, 5800, 45747.97916666666),
    @(97, 5740, 45747.98958333334),
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
}
